# Add the new "Dove Springs Neighborhood Public Health Facility" location
# as row 32 of the Austin_Public_Health_Locations sheet, then leave the
# selection on E32 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "Dove Springs Neighborhood Public Health Facility"
$ws.Range("B32").Value = "Neighborhood Center"
$ws.Range("C32").Value = "5106 Village Square Drive, Austin, Texas 78744"
$ws.Range("D32").Value = 78744
$ws.Range("E32").Value = 30.188171000000001
$ws.Range("F32").Value = -97.744249999999994

$ws.Range("E32").Select()
